# 1.5.4.xlsx update:
#  - remove the helper "Parameters" worksheet (its reference data is no
#    longer needed once the figures it backed are baked straight into
#    the indicator sheet)
#  - extend "Лист1" with four more reporting years (2020-2023) in
#    columns E:H, reusing the formatting already used for 2019 (col D)

$excel.DisplayAlerts = $false

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# --- drop the now-unused "Parameters" sheet ------------------------------
$params = $wb.Worksheets.Item("Parameters")
$params.Delete()

# --- copy the formatting of column D (2019) across to E:H ---------------
$ws.Range("D3:D6").Copy() | Out-Null
$ws.Range("E3:H6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- year headers ---------------------------------------------------------
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# --- row 4: number of local governments -----------------------------------
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# --- row 5: proportion of local governments (%) ----------------------------
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# --- row 6: number of local governments with DRR strategies ---------------
$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# restore the selection Excel leaves behind after this kind of edit
$ws.Range("D9").Select()
